$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds plain-text numbers (e.g. "26.774.60", "0.505").
# Force text format first so Excel does not reinterpret them as real numbers
# when we assign the new values below.
$ws.Range("D2:D47").NumberFormat = "@"
$ws.Range("B48:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.774.60"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "1.647.82"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").Value = "  +0.60%  "
$ws.Range("D5").Value = "216.67"
$ws.Range("E5").Value = "  +0.76%  "
$ws.Range("D6").Value = "0.505"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("E7").Value = "  +0.56%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -0.40%  "
$ws.Range("D10").Value = "19.27"
$ws.Range("E10").Value = "  -0.30%  "
$ws.Range("D11").Value = "0.0844"
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("D12").Value = "1.874.15"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").Value = "1.646.04"
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("D14").Value = "4.21"
$ws.Range("D15").Value = "0.532"
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").Value = "65.67"
$ws.Range("E16").Value = "  -0.49%  "
$ws.Range("D17").Value = "26.781.13"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("E18").Value = "  -0.29%  "
$ws.Range("D19").Value = "217.35"
$ws.Range("E19").Value = "  -0.73%  "
$ws.Range("E20").Value = "  +0.64%  "
$ws.Range("D21").Value = "4.37"
$ws.Range("E21").Value = "  +0.26%  "
$ws.Range("D22").Value = "2.42"
$ws.Range("E22").Value = "  +15.73%  "
$ws.Range("E23").Value = "  -0.69%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "145.67"
$ws.Range("E25").Value = "  -1.44%  "
$ws.Range("E26").Value = "  +0.71%  "
$ws.Range("D28").Value = "7.21"
$ws.Range("E28").Value = "  +3.76%  "
$ws.Range("D29").Value = "15.84"
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  +0.67%  "
$ws.Range("E31").Value = "  +0.88%  "
$ws.Range("D32").Value = "3.35"
$ws.Range("E32").Value = "  -0.80%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("E34").Value = "  +1.17%  "
$ws.Range("D35").Value = "1.276.47"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("D36").Value = "2.43"
$ws.Range("E36").Value = "  +2.07%  "
$ws.Range("E37").Value = "  +0.45%  "
$ws.Range("E38").Value = "  +5.23%  "
$ws.Range("D39").Value = "0.833"
$ws.Range("E39").Value = "  +3.06%  "
$ws.Range("E40").Value = "  +0.62%  "
$ws.Range("D41").Value = "0.818"
$ws.Range("E41").Value = "  +1.78%  "
$ws.Range("E42").Value = "  -1.34%  "
$ws.Range("E43").Value = "  +1.19%  "
$ws.Range("D44").Value = "1.799.32"
$ws.Range("E44").Value = "  +0.64%  "
$ws.Range("D45").Value = "92.12"
$ws.Range("E45").Value = "  -1.78%  "
$ws.Range("D46").Value = "59.62"
$ws.Range("E46").Value = "  +6.46%  "
$ws.Range("D47").Value = "1.62"
$ws.Range("E47").Value = "  +0.99%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0104"
$ws.Range("E48").Value = "  +0.49%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.0516"
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "7.76"
$ws.Range("E50").Value = "  +1.23%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.0983"
$ws.Range("E51").Value = "  +1.56%  "

# Re-apply the default (unstyled) look now that the text values are locked in,
# matching the original workbook formatting.
$ws.Range("D2:D47").Style = "Normal"
$ws.Range("B48:D51").Style = "Normal"
